$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("I2").Value = 0.1170416182565993
$ws.Range("J2").Value = 0.1170416182565993
$ws.Range("M2").Value = 6.174188000000001
$ws.Range("N2").Value = 18.522564
$ws.Range("O2").Value = 0.2521113718502555
$ws.Range("P2").Value = 0.2521113718502555
$ws.Range("Q2").Value = 1.600135491082667
$ws.Range("R2").Value = 14.401219419744
$ws.Range("S2").Value = 0.02950752294224516
$ws.Range("T2").Value = 0.02950752294224516
$ws.Range("I3").Value = 0.1170416182565993
$ws.Range("J3").Value = 0.1170416182565993
$ws.Range("M3").Value = 8.118224666666666
$ws.Range("O3").Value = 0.331492458231255
$ws.Range("P3").Value = 0.3314924582312551
$ws.Range("Q3").Value = 2.103962401811555
$ws.Range("S3").Value = 0.03879841375124424
$ws.Range("T3").Value = 0.03879841375124424
$ws.Range("I4").Value = 0.1170416182565993
$ws.Range("J4").Value = 0.1170416182565993
$ws.Range("M4").Value = 4.957885
$ws.Range("N4").Value = 14.873655
$ws.Range("O4").Value = 0.2024459230632115
$ws.Range("P4").Value = 0.2024459230632115
$ws.Range("Q4").Value = 1.284911918653333
$ws.Range("R4").Value = 11.56420726788
$ws.Range("S4").Value = 0.02369459844476927
$ws.Range("T4").Value = 0.02369459844476928
$ws.Range("I5").Value = 0.1170416182565993
$ws.Range("J5").Value = 0.1170416182565993
$ws.Range("M5").Value = 1.245063666666667
$ws.Range("N5").Value = 3.735191
$ws.Range("O5").Value = 0.05083983659782346
$ws.Range("P5").Value = 0.05083983659782347
$ws.Range("Q5").Value = 0.3226773401928889
$ws.Range("R5").Value = 2.904096061736
$ws.Range("S5").Value = 0.00595037674731034
$ws.Range("T5").Value = 0.005950376747310341
$ws.Range("I6").Value = 0.1170416182565993
$ws.Range("J6").Value = 0.1170416182565993
$ws.Range("M6").Value = 3.994561333333333
$ws.Range("N6").Value = 11.983684
$ws.Range("O6").Value = 0.1631104102574545
$ws.Range("P6").Value = 0.1631104102574545
$ws.Range("Q6").Value = 1.035251819473778
$ws.Range("R6").Value = 9.317266375263999
$ws.Range("S6").Value = 0.01909070637103028
$ws.Range("T6").Value = 0.01909070637103029
$ws.Range("G7").Value = 1.708219666666666
$ws.Range("H7").Value = 5.124658999999999
$ws.Range("I7").Value = 0.7714488336573383
$ws.Range("J7").Value = 0.7714488336573383
$ws.Range("M7").Value = 6.174188000000001
$ws.Range("N7").Value = 18.522564
$ws.Range("O7").Value = 0.2521113718502555
$ws.Range("P7").Value = 0.2521113718502555
$ws.Range("Q7").Value = 10.54686936729733
$ws.Range("R7").Value = 94.921824305676
$ws.Range("S7").Value = 0.1944910237656311
$ws.Range("T7").Value = 0.1944910237656311
$ws.Range("G8").Value = 1.708219666666666
$ws.Range("H8").Value = 5.124658999999999
$ws.Range("I8").Value = 0.7714488336573383
$ws.Range("J8").Value = 0.7714488336573383
$ws.Range("M8").Value = 8.118224666666666
$ws.Range("O8").Value = 0.331492458231255
$ws.Range("P8").Value = 0.3314924582312551
$ws.Range("Q8").Value = 13.86771103401844
$ws.Range("R8").Value = 124.809399306166
$ws.Range("S8").Value = 0.2557294702687056
$ws.Range("T8").Value = 0.2557294702687056
$ws.Range("G9").Value = 1.708219666666666
$ws.Range("H9").Value = 5.124658999999999
$ws.Range("I9").Value = 0.7714488336573383
$ws.Range("J9").Value = 0.7714488336573383
$ws.Range("M9").Value = 4.957885
$ws.Range("N9").Value = 14.873655
$ws.Range("O9").Value = 0.2024459230632115
$ws.Range("P9").Value = 0.2024459230632115
$ws.Range("Q9").Value = 8.469156662071667
$ws.Range("R9").Value = 76.22240995864499
$ws.Range("S9").Value = 0.1561766712257978
$ws.Range("T9").Value = 0.1561766712257978
$ws.Range("G10").Value = 1.708219666666666
$ws.Range("H10").Value = 5.124658999999999
$ws.Range("I10").Value = 0.7714488336573383
$ws.Range("J10").Value = 0.7714488336573383
$ws.Range("M10").Value = 1.245063666666667
$ws.Range("N10").Value = 3.735191
$ws.Range("O10").Value = 0.05083983659782346
$ws.Range("P10").Value = 0.05083983659782347
$ws.Range("Q10").Value = 2.126842241652111
$ws.Range("R10").Value = 19.141580174869
$ws.Range("S10").Value = 0.03922033264672057
$ws.Range("T10").Value = 0.03922033264672058
$ws.Range("G11").Value = 1.708219666666666
$ws.Range("H11").Value = 5.124658999999999
$ws.Range("I11").Value = 0.7714488336573383
$ws.Range("J11").Value = 0.7714488336573383
$ws.Range("M11").Value = 3.994561333333333
$ws.Range("N11").Value = 11.983684
$ws.Range("O11").Value = 0.1631104102574545
$ws.Range("P11").Value = 0.1631104102574545
$ws.Range("Q11").Value = 6.823588229306222
$ws.Range("R11").Value = 61.41229406375599
$ws.Range("S11").Value = 0.1258313357504832
$ws.Range("T11").Value = 0.1258313357504832
$ws.Range("G12").Value = 0.1229426666666667
$ws.Range("H12").Value = 0.368828
$ws.Range("I12").Value = 0.05552211970009493
$ws.Range("J12").Value = 0.05552211970009493
$ws.Range("M12").Value = 6.174188000000001
$ws.Range("N12").Value = 18.522564
$ws.Range("O12").Value = 0.2521113718502555
$ws.Range("P12").Value = 0.2521113718502555
$ws.Range("Q12").Value = 0.7590711372213333
$ws.Range("R12").Value = 6.831640234992001
$ws.Range("S12").Value = 0.01399775776562503
$ws.Range("T12").Value = 0.01399775776562503
$ws.Range("G13").Value = 0.1229426666666667
$ws.Range("H13").Value = 0.368828
$ws.Range("I13").Value = 0.05552211970009493
$ws.Range("J13").Value = 0.05552211970009493
$ws.Range("M13").Value = 8.118224666666666
$ws.Range("O13").Value = 0.331492458231255
$ws.Range("P13").Value = 0.3314924582312551
$ws.Range("Q13").Value = 0.998076189119111
$ws.Range("R13").Value = 8.982685702071999
$ws.Range("S13").Value = 0.01840516394559446
$ws.Range("T13").Value = 0.01840516394559446
$ws.Range("G14").Value = 0.1229426666666667
$ws.Range("H14").Value = 0.368828
$ws.Range("I14").Value = 0.05552211970009493
$ws.Range("J14").Value = 0.05552211970009493
$ws.Range("M14").Value = 4.957885
$ws.Range("N14").Value = 14.873655
$ws.Range("O14").Value = 0.2024459230632115
$ws.Range("P14").Value = 0.2024459230632115
$ws.Range("Q14").Value = 0.6095356029266666
$ws.Range("R14").Value = 5.485820426339999
$ws.Range("S14").Value = 0.01124022677311184
$ws.Range("T14").Value = 0.01124022677311184
$ws.Range("G15").Value = 0.1229426666666667
$ws.Range("H15").Value = 0.368828
$ws.Range("I15").Value = 0.05552211970009493
$ws.Range("J15").Value = 0.05552211970009493
$ws.Range("M15").Value = 1.245063666666667
$ws.Range("N15").Value = 3.735191
$ws.Range("O15").Value = 0.05083983659782346
$ws.Range("P15").Value = 0.05083983659782347
$ws.Range("Q15").Value = 0.1530714473497778
$ws.Range("R15").Value = 1.377643026148
$ws.Range("S15").Value = 0.002822735493117621
$ws.Range("T15").Value = 0.002822735493117622
$ws.Range("G16").Value = 0.1229426666666667
$ws.Range("H16").Value = 0.368828
$ws.Range("I16").Value = 0.05552211970009493
$ws.Range("J16").Value = 0.05552211970009493
$ws.Range("M16").Value = 3.994561333333333
$ws.Range("N16").Value = 11.983684
$ws.Range("O16").Value = 0.1631104102574545
$ws.Range("P16").Value = 0.1631104102574545
$ws.Range("Q16").Value = 0.4911020224835556
$ws.Range("R16").Value = 4.419918202352
$ws.Range("S16").Value = 0.009056235722645979
$ws.Range("T16").Value = 0.009056235722645979
$ws.Range("G17").Value = 0.123973
$ws.Range("H17").Value = 0.371919
$ws.Range("I17").Value = 0.05598742838596747
$ws.Range("J17").Value = 0.05598742838596747
$ws.Range("M17").Value = 6.174188000000001
$ws.Range("N17").Value = 18.522564
$ws.Range("O17").Value = 0.2521113718502555
$ws.Range("P17").Value = 0.2521113718502555
$ws.Range("Q17").Value = 0.7654326089240001
$ws.Range("R17").Value = 6.888893480316001
$ws.Range("S17").Value = 0.01411506737675419
$ws.Range("T17").Value = 0.01411506737675419
$ws.Range("G18").Value = 0.123973
$ws.Range("H18").Value = 0.371919
$ws.Range("I18").Value = 0.05598742838596747
$ws.Range("J18").Value = 0.05598742838596747
$ws.Range("M18").Value = 8.118224666666666
$ws.Range("O18").Value = 0.331492458231255
$ws.Range("P18").Value = 0.3314924582312551
$ws.Range("Q18").Value = 1.006440666600667
$ws.Range("R18").Value = 9.057965999405999
$ws.Range("S18").Value = 0.0185594102657107
$ws.Range("T18").Value = 0.01855941026571071
$ws.Range("G19").Value = 0.123973
$ws.Range("H19").Value = 0.371919
$ws.Range("I19").Value = 0.05598742838596747
$ws.Range("J19").Value = 0.05598742838596747
$ws.Range("M19").Value = 4.957885
$ws.Range("N19").Value = 14.873655
$ws.Range("O19").Value = 0.2024459230632115
$ws.Range("P19").Value = 0.2024459230632115
$ws.Range("Q19").Value = 0.614643877105
$ws.Range("R19").Value = 5.531794893944999
$ws.Range("S19").Value = 0.01133442661953263
$ws.Range("T19").Value = 0.01133442661953263
$ws.Range("G20").Value = 0.123973
$ws.Range("H20").Value = 0.371919
$ws.Range("I20").Value = 0.05598742838596747
$ws.Range("J20").Value = 0.05598742838596747
$ws.Range("M20").Value = 1.245063666666667
$ws.Range("N20").Value = 3.735191
$ws.Range("O20").Value = 0.05083983659782346
$ws.Range("P20").Value = 0.05083983659782347
$ws.Range("Q20").Value = 0.1543542779476667
$ws.Range("R20").Value = 1.389188501529
$ws.Range("S20").Value = 0.002846391710674929
$ws.Range("T20").Value = 0.002846391710674929
$ws.Range("G21").Value = 0.123973
$ws.Range("H21").Value = 0.371919
$ws.Range("I21").Value = 0.05598742838596747
$ws.Range("J21").Value = 0.05598742838596747
$ws.Range("M21").Value = 3.994561333333333
$ws.Range("N21").Value = 11.983684
$ws.Range("O21").Value = 0.1631104102574545
$ws.Range("P21").Value = 0.1631104102574545
$ws.Range("Q21").Value = 0.4952177521773333
$ws.Range("R21").Value = 4.456959769596
$ws.Range("S21").Value = 0.009132132413295005
$ws.Range("T21").Value = 0.009132132413295005
